$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '63.013.07'
Set-TextValue 'E2' '  -1.99%  '
Set-TextValue 'D3' '2.573.01'
Set-TextValue 'E3' '  -3.05%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '586.07'
Set-TextValue 'E5' '  -3.31%  '
Set-TextValue 'D6' '147.60'
Set-TextValue 'E6' '  -3.54%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.583'
Set-TextValue 'E8' '  -1.31%  '
Set-TextValue 'E9' '  -1.23%  '
Set-TextValue 'E10' '  +1.43%  '
Set-TextValue 'D11' '0.379'
Set-TextValue 'E11' '  -1.83%  '
Set-TextValue 'E12' '  -0.85%  '
Set-TextValue 'D13' '27.20'
Set-TextValue 'E13' '  -3.53%  '
Set-TextValue 'D14' '3.032.79'
Set-TextValue 'E14' '  -3.12%  '
Set-TextValue 'D15' '62.902.83'
Set-TextValue 'E16' '  +2.22%  '
Set-TextValue 'D17' '2.562.05'
Set-TextValue 'E17' '  -3.25%  '
Set-TextValue 'E18' '  -0.62%  '
Set-TextValue 'D19' '4.63'
Set-TextValue 'E19' '  -0.12%  '
Set-TextValue 'D20' '342.50'
Set-TextValue 'E20' '  -1.60%  '
Set-TextValue 'E21' '  -2.40%  '
Set-TextValue 'E22' '  -0.18%  '
Set-TextValue 'D23' '66.56'
Set-TextValue 'E23' '  -0.18%  '
Set-TextValue 'E24' '  -3.80%  '
Set-TextValue 'D25' '9.03'
Set-TextValue 'E25' '  -3.83%  '
Set-TextValue 'E26' '  -4.74%  '
Set-TextValue 'D27' '550.05'
Set-TextValue 'D28' '7.98'
Set-TextValue 'E28' '  -2.48%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  +0.16%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D30' '0.161'
Set-TextValue 'E30' '  -2.23%  '
Set-TextValue 'E31' '  -2.44%  '
Set-TextValue 'D32' '0.0₃0841'
Set-TextValue 'E32' '  -2.39%  '
Set-TextValue 'E33' '  -2.84%  '
Set-TextValue 'D34' '5.12'
Set-TextValue 'E34' '  -4.29%  '
Set-TextValue 'D35' '165.20'
Set-TextValue 'E35' '  -1.86%  '
Set-TextValue 'D36' '0.409'
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'D37' '0.999'
Set-TextValue 'E37' '  -0.03%  '
Set-TextValue 'E38' '  -0.64%  '
Set-TextValue 'E39' '  -4.77%  '
Set-TextValue 'E40' '  -0.04%  '
Set-TextValue 'D41' '164.70'
Set-TextValue 'E41' '  -0.98%  '
Set-TextValue 'D42' '39.56'
Set-TextValue 'E42' '  -1.28%  '
Set-TextValue 'D43' '3.92'
Set-TextValue 'E43' '  +1.75%  '
Set-TextValue 'E44' '  +0.77%  '
Set-TextValue 'D45' '22.38'
Set-TextValue 'E45' '  +2.12%  '
Set-TextValue 'E46' '  -0.77%  '
Set-TextValue 'E47' '  +0.58%  '
Set-TextValue 'D48' '0.0245'
Set-TextValue 'E48' '  -0.74%  '
Set-TextValue 'E49' '  -0.92%  '
Set-TextValue 'D50' '18.80'
Set-TextValue 'E50' '  -1.37%  '
Set-TextValue 'E51' '  +10.17%  '
